$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 9260251.38139425

$ws.Range("B2:F7").Value = $newValue
